# "Generate Report for Handoff"
#
# The localization-status report gets a fresh handoff-report run for the
# last tracked file (f43d7ba1-db92-40a5-b276-ec4111773384): its "latest
# handoff" timestamp is refreshed on the Overview sheet as well as on each
# per-locale detail sheet (zh-cn, de-de). These are plain text cells (the
# source data isn't valid real dates - e.g. month "30") so they must stay
# strings, not get coerced into Excel date serials.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest Handoff Date" column (D), last data row (7) -
# the row for f43d7ba1-db92-40a5-b276-ec4111773384.md
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value = "2016-30-18 00:30:59"

# zh-cn detail sheet: "Latest Handoff Datetime" column (E), last data row (7)
$zhCn = $wb.Worksheets.Item("zh-cn")
$zhCn.Range("E7").Value = "2016-03-18 00:30:56"

# de-de detail sheet: "Latest Handoff Datetime" column (E), last data row (7)
$deDe = $wb.Worksheets.Item("de-de")
$deDe.Range("E7").Value = "2016-03-18 00:30:59"
